$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 60250
$ws.Range("B4").Value = 11181
$ws.Range("C4").Value = 126892
$ws.Range("D4").Value = 4.415023487107064
$ws.Range("E4").Value = 0.8193878517244749
$ws.Range("F4").Value = 9.29842255875332
$ws.Range("G4").Value = 0.2363384724149416
$ws.Range("H4").Value = 0.08512938799278082
$ws.Range("I4").Value = 0.3385016569214514
$ws.Range("J4").Value = 0.1415727584975331
$ws.Range("K4").Value = 0.0724565009258427
$ws.Range("L4").Value = 0.2527275740391631
$ws.Range("M4").Value = 125973
$ws.Range("N4").Value = 13263
$ws.Range("O4").Value = 194670
$ws.Range("P4").Value = 9.650647125870528
$ws.Range("Q4").Value = 1.01610534491522
$ws.Range("R4").Value = 14.91350116887804
$ws.Range("S4").Value = 0.2994429132643171
$ws.Range("T4").Value = 0.07833578680451396
$ws.Range("U4").Value = 0.4675508886010589
$ws.Range("V4").Value = 0.3231122866871788
$ws.Range("W4").Value = 0.08778803241880449
$ws.Range("X4").Value = 0.4891233426515225
